$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "US Gross Domestic Income and Equity Indices.xlsx",
    "Bitty, SPX & GM2 fitted trends.xlsx",
    ".DS_Store",
    "Monetary Expansion and Risk Asset Activity.xlsx",
    "Other deposit liabilities (ODL) shows where M2 is going.xlsx",
    " China domestic credit & M2.xlsx",
    "US financial conditions, employment & equities.xlsx",
    "GDP vs GDI.xlsx",
    "Financial conditions & employment, US.xlsx",
    "Bank credit and M2 U.S.xlsx",
    "Bank credit growth, inflation and M2 U.S.A.xlsx",
    "Total construction spending in the U.S, nominal and real.xlsx",
    "Bitcoin price change is fueled by global monetary growth.xlsx",
    "Gross domestic product (GDP) vs gross domestic income (GDI), yearly growth..xlsx",
    "U.S GDP and GDI.xlsx",
    "Bitcoin price change is fueled by global monetary growth (with forecast).xlsx",
    "Personal savings U.S (BEA).xlsx",
    "Real Incomes U.S (BEA).xlsx",
    "US and Global Monetary Aggregates.xlsx",
    "US Equity Indices.xlsx"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
